$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 - this shifts the existing rows 2-104 down to
# 3-105, carrying their content/height/style along (matches the diff's row
# shift pattern exactly).
$ws.Rows.Item(2).Insert()

# Copy the format of what is now row 9 (the original "Day 7" row, formerly
# row 8) onto the new row 2, so A2:E2 pick up styles s11/s12/s13/s4/s4.
$ws.Range("A9:E9").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(2).RowHeight = 88.5

# Populate the new row 2 with the "Day 7" entry, now marked completed with a
# remark (the task/date text duplicate what's already in row 9).
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "20.05.2025"
$ws.Range("C2").Value = "React - Introduction to React , Difference Class and Functional Component"
$ws.Range("D2").Value = "Completed"
$ws.Range("E2").Value = "Learnt basics of react and different components"

# Update the view: selection moves to E2 and the scrolled-to top-left cell
# resets (matches the diff removing topLeftCell="A7").
$ws.Range("E2").Select()
